# Add the new "2022-Q3" quarterly sheet and update the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet right after "总计" (i.e. before the existing
#    "2022-Q2" sheet) and name it "2022-Q3". Duplicating the "2022-Q2" sheet
#    (rather than adding a blank one) keeps all of its sheet-level
#    formatting (outline props, page margins, header/index-column styles)
#    intact, matching the look of the other quarterly sheets.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$beforeSheet.Copy($beforeSheet)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The duplicated sheet has 16 data rows (same as "2022-Q2"); 2022-Q3 only
# needs 5, so drop the extra rows before filling in the real values.
$q3.Range("A7:H17").Clear()

# ---------------------------------------------------------------------------
# 2. Populate the header row.
# ---------------------------------------------------------------------------
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3. Populate the fund holding rows for 2022-Q3.
# ---------------------------------------------------------------------------
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'001305"
$q3.Range("C2").Value = "九泰天富改革新动力混合A"
$q3.Range("D2").Value = "'1.86"
$q3.Range("E2").Value = "'94.71"
$q3.Range("F2").Value = "'8.70"
$q3.Range("G2").Value = "'0.1618"
$q3.Range("H2").Value = 3

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'001844"
$q3.Range("C3").Value = "九泰久益灵活配置混合C"
$q3.Range("D3").Value = "'0.98"
$q3.Range("E3").Value = "'93.32"
$q3.Range("F3").Value = "'8.03"
$q3.Range("G3").Value = "'0.0787"
$q3.Range("H3").Value = 5

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'001782"
$q3.Range("C4").Value = "九泰久益灵活配置混合A"
$q3.Range("D4").Value = "'0.53"
$q3.Range("E4").Value = "'93.32"
$q3.Range("F4").Value = "'8.03"
$q3.Range("G4").Value = "'0.0426"
$q3.Range("H4").Value = 5

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "'004332"
$q3.Range("C5").Value = "恒生前海沪港深新兴产业精选混合"
$q3.Range("D5").Value = "'0.47"
$q3.Range("E5").Value = "'92.74"
$q3.Range("F5").Value = "'8.28"
$q3.Range("G5").Value = "'0.0389"
$q3.Range("H5").Value = 4

$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "'009912"
$q3.Range("C6").Value = "九泰天富改革新动力混合C"
$q3.Range("D6").Value = "'0.17"
$q3.Range("E6").Value = "'94.71"
$q3.Range("F6").Value = "'8.70"
$q3.Range("G6").Value = "'0.0148"
$q3.Range("H6").Value = 3

# ---------------------------------------------------------------------------
# 4. Update the "总计" summary sheet: insert a new row for 2022-Q3 above the
#    existing 2022-Q2 row and push everything else down one row, then
#    re-number the index column (A) sequentially.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The inserted row inherits formatting from the header row above; reset it
# to match the plain (unstyled) data rows by copying the format from the
# row directly beneath (the former 2022-Q2 row, now row 3).
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.34

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7

# ---------------------------------------------------------------------------
# 5. Restore the originally-active sheet/selection so the workbook's view
#    state isn't left pointing at a sheet we only touched incidentally.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
